$d = $word.ActiveDocument

# --- Change 1 ---
# Merge the two runs split by a lastRenderedPageBreak in the "Cần một
# phương pháp ..." paragraph back into a single run of text by replacing
# the text that spans the break (the trailing space of the first run plus
# the start of the second run) with itself. Word's Find/Replace collapses
# the runs either side of the match, eliminating the now-redundant
# lastRenderedPageBreak run boundary and the trailing "preserve" space run.
$d.Content.Find.Execute("bản đồ quét sang vector", $true, $false, $false, $false, $false, $true, 1, $false, "bản đồ quét sang vector", 2) | Out-Null

# --- Change 2 ---
# Same pattern for the "Toàn bộ tập hợp 256 trạng thái xung quanh ..."
# paragraph.
$d.Content.Find.Execute("được gọi là nhóm mẫu", $true, $false, $false, $false, $false, $true, 1, $false, "được gọi là nhóm mẫu", 2) | Out-Null

# --- Change 3 ---
# Same merge pattern for the "Mục đích của hoạt động lấp đầy ..."
# paragraph.
$d.Content.Find.Execute("trong các đường nét ban đầu", $true, $false, $false, $false, $false, $true, 1, $false, "trong các đường nét ban đầu", 2) | Out-Null

# In addition, this paragraph originally continued with a manual line
# break run (<w:br/>) right after the merged text. The edit splits that
# off into its own new paragraph. Locate the manual line break
# (character code 11) and insert a paragraph break immediately before it.
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $idx = $t.IndexOf([char]11)
    if ($idx -ge 0) {
        $absPos = $p.Range.Start + $idx
        $ins = $d.Range($absPos, $absPos)
        $ins.InsertParagraphBefore()
        break
    }
}
